# "all stats debug in progress"
# Populate the per-phoneme-type percentage rows (5 and 6) on the "Vowels"
# sheet, rescale row 4 from raw counts to percentages, and fill in the
# missing counts on row 5/6 of the "Cons manner" sheet (and fix row 4's
# stray value there too).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Vowels" (sheet1): columns D..K
# ---------------------------------------------------------------------
$wsVowels = $wb.Worksheets.Item("Vowels")

$vowelCols = @("D", "E", "F", "G", "H", "I", "J", "K")

$row4Values = @(0.4, 0.03333333333333333, 0.03333333333333333, 0.5, 0.4, 0.7666666666666667, 0.06666666666666667, 0.43333333333333335)
$row5Values = @(0.11333333333333333, 0.006666666666666667, 0.013333333333333334, 0.16, 0.13333333333333333, 0.26, 0.02666666666666667, 0.14)
$row6Values = @(0.5666666666666667, 0.03333333333333333, 0.06666666666666667, 0.8, 0.6666666666666666, 1.3, 0.13333333333333333, 0.7)

for ($i = 0; $i -lt $vowelCols.Length; $i++) {
    $col = $vowelCols[$i]

    $c4 = $wsVowels.Range($col + "4")
    $c4.Value = $row4Values[$i]
    $c4.NumberFormat = "0.0%"

    $c5 = $wsVowels.Range($col + "5")
    $c5.Value = $row5Values[$i]
    $c5.NumberFormat = "0.0%"

    $c6 = $wsVowels.Range($col + "6")
    $c6.Value = $row6Values[$i]
    $c6.NumberFormat = "0.0"
}

# ---------------------------------------------------------------------
# Sheet "Cons manner" (sheet2): column E
# ---------------------------------------------------------------------
$wsManner = $wb.Worksheets.Item("Cons manner")

$wsManner.Range("E4").Value = 23.0
$wsManner.Range("E5").Value = 36.0
$wsManner.Range("E6").Value = 36.0
